$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G14").Value = "Acierto"
$ws.Range("H14").Value = 1.63
